$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.107.96"
$ws.Range("E2").Value = "  +0.37%  "

# Row 3
$ws.Range("D3").Value = "2.535.77"
$ws.Range("E3").Value = "  +0.16%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.21"
$ws.Range("E5").Value = "  -0.37%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.84"
$ws.Range("E6").Value = "  -1.55%  "

# Row 7
$ws.Range("E7").Value = "  +0.10%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.524"
$ws.Range("E8").Value = "  -1.11%  "

# Row 9
$ws.Range("D9").Value = "2.538.35"
$ws.Range("E9").Value = "  +0.27%  "

# Row 10
$ws.Range("E10").Value = "  -3.31%  "

# Row 11
$ws.Range("E11").Value = "  +1.64%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.344"
$ws.Range("E12").Value = "  -0.17%  "

# Row 13
$ws.Range("E13").Value = "  -2.78%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.46"
$ws.Range("E14").Value = "  -1.56%  "

# Row 15
$ws.Range("E15").Value = "  +0.03%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000175"
$ws.Range("E16").Value = "  -1.57%  "

# Row 17
$ws.Range("D17").Value = "68.258.48"
$ws.Range("E17").Value = "  +0.91%  "

# Row 18
$ws.Range("D18").Value = "2.541.10"
$ws.Range("E18").Value = "  +0.08%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.91"
$ws.Range("E19").Value = "  +3.84%  "

# Row 20
$ws.Range("E20").Value = "  -0.35%  "

# Row 21
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "365.42"
$ws.Range("E21").Value = "  +1.32%  "

# Row 22
$ws.Range("B22").Value = "Polkadot"
$ws.Range("C22").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.16"
$ws.Range("E22").Value = "  -1.19%  "

# Row 23
$ws.Range("B23").Value = "Binance-PegBSC-USD"
$ws.Range("C23").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.35"
$ws.Range("E23").Value = "  +35.69%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.53"
$ws.Range("E24").Value = "  -2.70%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.04"
$ws.Range("E25").Value = "  +1.57%  "

# Row 26
$ws.Range("E26").Value = "  +0.00%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.88"
$ws.Range("E27").Value = "  -5.15%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.87"
$ws.Range("E28").Value = "  -4.47%  "

# Row 29
$ws.Range("D29").Value = "2.672.81"
$ws.Range("E29").Value = "  +0.70%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0956"
$ws.Range("E30").Value = "  -3.40%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "532.87"
$ws.Range("E31").Value = "  -4.01%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.26"
$ws.Range("E32").Value = "  -0.06%  "

# Row 33
$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.86"
$ws.Range("E33").Value = "  -0.14%  "

# Row 34
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.30"
$ws.Range("E34").Value = "  -4.08%  "

# Row 35
$ws.Range("E35").Value = "  -1.78%  "

# Row 36
$ws.Range("E36").Value = "  +0.11%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.85"
$ws.Range("E37").Value = "  +2.60%  "

# Row 38
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.28"
$ws.Range("E38").Value = "  +2.64%  "

# Row 39
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.44"
$ws.Range("E39").Value = "  -2.82%  "

# Row 40
$ws.Range("E40").Value = "  +0.20%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.78"
$ws.Range("E41").Value = "  -1.84%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.10"
$ws.Range("E42").Value = "  -1.67%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.347"
$ws.Range("E43").Value = "  -2.77%  "

# Row 44
$ws.Range("E44").Value = "  +0.19%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.48"
$ws.Range("E45").Value = "  -1.98%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.41"
$ws.Range("E46").Value = "  -1.36%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.51"
$ws.Range("E47").Value = "  +0.77%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.553"
$ws.Range("E48").Value = "  -1.72%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.70"
$ws.Range("E49").Value = "  -0.73%  "

# Row 50
$ws.Range("D50").Value = "0.0₆0275"
$ws.Range("E50").Value = "  -2.02%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.71"
$ws.Range("E51").Value = "  +0.95%  "
